$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Fitness) values for rows 2-8 to 3986
$ws.Range("C2:C8").Value = 3986

# Update column C (Fitness) values for rows 9-12 to 4046
$ws.Range("C9:C12").Value = 4046
